# Update excess mortality week 44
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The totals row currently lives on row 37, with rows 37-39 empty after it.
# Insert 3 blank rows above the totals row so it moves down to row 40,
# leaving row 36 free for the new week 44 data.
$ws.Rows("37:39").Insert()

# Fill in the new week 44 data on row 36 (continues the existing table).
$ws.Range("F36").Value = 44
$ws.Range("G36").Value = 3617
$ws.Range("H36").Value = 2889
$ws.Range("I36").Formula = "=G36-H36"

# Move the view so the newly added row is visible, matching the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("I36").Select()
